$wb = $excel.ActiveWorkbook

# Fix typo in shared string used by PSA_LOLO!B1 ("psa_loll_40" -> "psa_lolo_40")
$psaSheet = $wb.Worksheets.Item("PSA_LOLO")
$psaSheet.Range("B1").Value = "psa_lolo_40"

# Make PSA_LOLO the active (selected) sheet/tab instead of OverallRebateEfficiency
$psaSheet.Activate()
